$wb = $excel.ActiveWorkbook

$wsErl   = $wb.Worksheets.Item("Erläuterung")
$ws0301  = $wb.Worksheets.Item("03.01.21")
$wsTag   = $wb.Worksheets.Item("Impfungen_proTag")

# --- "03.01.21" sheet: updated figures for Rheinland-Pfalz (row 12) ---
$ws0301.Range("B12").Value = 8282
$ws0301.Range("C12").Value = 658
$ws0301.Range("D12").Value = 2.0230083614584911
$ws0301.Range("F12").Value = 4213
$ws0301.Range("H12").Value = 4069

# --- "03.01.21" sheet: updated totals (row 18 "Gesamt") ---
$ws0301.Range("B18").Value = 265986
$ws0301.Range("C18").Value = 22892
$ws0301.Range("D18").Value = 3.1982267520474625
$ws0301.Range("F18").Value = 123103
$ws0301.Range("H18").Value = 114654

# --- "Impfungen_proTag" sheet: B9 becomes a formula pulling the new total ---
$wsTag.Range("B9").Formula = "='03.01.21'!C18"

# --- Selection / active-cell bookkeeping to match the saved workbook state ---
$wsTag.Range("B10").Select() | Out-Null

$wsErl.Range("K18").Select() | Out-Null

# Make "03.01.21" the active sheet/tab, with its new active cell selection
$ws0301.Activate() | Out-Null
$ws0301.Range("K14").Select() | Out-Null
